$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set row 2 (A2:H2) values to 0
$ws.Range("A2:H2").Value = 0

# Set row 13 (A13:H13) values to 0
$ws.Range("A13:H13").Value = 0

# Update the selection to A2:H2 with active cell A2
$ws.Range("A2:H2").Select()
